$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows below the existing data (row 30 was the last row),
# copying number/alignment formatting down from the row above - matches
# the existing table's per-column styles (D: applyFill style, I: left
# aligned boolean style).
$ws.Rows(31).Insert(-4121, 0)
$ws.Rows(32).Insert(-4121, 0)

# Entered first: John Doe / john.doe@xyz.com (ends up placed in row 32)
$ws.Range("A32").Value = 110031
$ws.Range("B32").Value = 9317596767
$ws.Range("C32").Value = "John Doe"
$ws.Range("D32").Value = "john.doe@xyz.com"
$ws.Range("E32").Value = 818876431
$ws.Range("F32").Value = "ACT"
$ws.Range("G32").Value = "eng"
$ws.Range("H32").Value = "PWD"
$ws.Range("I32").Value = $true
$ws.Range("J32").Value = "superadmin"
$ws.Range("K32").Value = "now()"
$ws.Range("L32").Value = "now()"

# Entered second: Jane Smith / jane.smith@xyz.com (ends up placed in row 31)
$ws.Range("A31").Value = 110030
$ws.Range("B31").Value = 9317596768
$ws.Range("C31").Value = "Jane Smith"
$ws.Range("D31").Value = "jane.smith@xyz.com"
$ws.Range("E31").Value = 818876432
$ws.Range("F31").Value = "ACT"
$ws.Range("G31").Value = "eng"
$ws.Range("H31").Value = "PWD"
$ws.Range("I31").Value = $true
$ws.Range("J31").Value = "superadmin"
$ws.Range("K31").Value = "now()"
$ws.Range("L31").Value = "now()"

$ws.Range("F30").Select()
